$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clean up OCR artifacts / typos in the Point Name column (column B) ---
$ws.Range("B2").Value  = 'RTU MONITOR JUMPER CONTROL INHIBIT'
$ws.Range("B3").Value  = 'STINGER SWITCH 115KV'
$ws.Range("B4").Value  = 'IINYO 115KV CB'
$ws.Range("B5").Value  = 'INYO 115KV CB'
$ws.Range("B7").Value  = 'HAIWEE-IINYOKERN 115KV CB'
$ws.Range("B8").Value  = 'COSO-HAIWEE-IINYOKERN 115KV CB'
$ws.Range("B10").Value = 'OXBOW 115KV CB'
$ws.Range("B11").Value = 'OXBOW CUST CB'
$ws.Range("B12").Value = 'DIXIE VALLEY CUST CB'
$ws.Range("B13").Value = 'NO. BANK 115KV CB'
$ws.Range("B15").Value = '115KV CB LOW GAS/AIR'
$ws.Range("B17").Value = 'TRANS BK SUDDENTPRESS RESET “1'
$ws.Range("B19").Value = 'NO. BANK C30 FAIL'
$ws.Range("B26").Value = '115KV_UNDERFREQ POT FAIL'
$ws.Range("B32").Value = '115KV DFR FAIL 1B19'
$ws.Range("B33").Value = '115KV DFR OPERATION 1820'
$ws.Range("B34").Value = '115KV UNDERFREQ RELAY PICK UP'
$ws.Range("B35").Value = '115KV UNDERFREQ RELAY FAIL'
$ws.Range("B38").Value = 'CASA D-SHERWIN REACTOR CS'
$ws.Range("B57").Value = 'TS5KV CB 11" LBFB RELAY PS FAIL'
$ws.Range("B72").Value = 'COSO HWE TGT'
$ws.Range("B73").Value = 'COSO HWE TGT'
$ws.Range("B74").Value = 'COSO HWE GRD TGT'

# --- Remove trailing rows that no longer apply (filtered out of the pointlist) ---
$ws.Rows("76:80").Delete()

# --- Narrow column B now that the longest OCR-garbled strings are gone ---
# (closest value achievable through the ColumnWidth/pixel-quantised API to
# the canonical target width of 68.567768 character-units)
$ws.Columns("B").ColumnWidth = 67.65
